$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows after row 10 (rows 11-13), pushing existing rows 11-14 down to 14-17
$ws.Range("A11:C13").Insert()

# Fill in the new rows with the new election entries, in the same order the
# original author entered them (matches the shared-strings append order)
$ws.Range("A12").Value = "Gobernatura 19"
$ws.Range("B12").Value = "gb_19"

# Fix the color value for "Senado 18" row (add missing leading '#')
$ws.Range("C10").Value = "#348cae4"

$ws.Range("A13").Value = "Distrito local 19"
$ws.Range("B13").Value = "dl_19"

$ws.Range("C12").Value = "#5a189a"
$ws.Range("C13").Value = "#6b9080"

$ws.Range("A11").Value = "Presidencia Municipal 19"
$ws.Range("B11").Value = "pm_19"
$ws.Range("C11").Value = "#d68c45"

# Update the selected cell to match the saved view state
$ws.Range("C11").Select()
